$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TitleBlockData")

# Helper: force a numeric-looking string ("10", "1", "4", ...) to be
# stored as literal TEXT rather than being auto-coerced into a number by
# Excel's normal type sniffing on Range.Value assignment. We briefly mark
# the cell as Text-formatted, write the value, then clear the formatting
# change and re-apply the sheet's normal "label" alignment (left/center/
# indent 1) so the cell lands back on the same style as its neighbours.
function Set-TextValue($rng, $text) {
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4108
    $rng.IndentLevel = 1
}

# New rows 15-18 need the same cell styling (left/center, indent 1) as the
# rest of the table body. Copy formats down from an already-styled block
# of rows (2-5) matched in size, so the paste doesn't tile beyond 4 rows.
$ws.Range("A2:E5").Copy()
$ws.Range("A15:E18").PasteSpecial(-4122)

# Re-populate the title-block property list in its new (sorted, extended)
# order. Existing rows shift down and four new properties are inserted
# (Angle_Units, Length_Units, Mass_Units, Number of sheets).
$ws.Range("A3").Value = 'Angle_Units'
$ws.Range("B3").Value = '°'
$ws.Range("C3").Value = ''
$ws.Range("D3").Value = ''
$ws.Range("E3").Value = ''

$ws.Range("A4").Value = 'DN'
$ws.Range("B4").Value = 'DN'
$ws.Range("C4").Value = ''
$ws.Range("D4").Value = ''
$ws.Range("E4").Value = ''

$ws.Range("A5").Value = 'DRAWING_TITLE'
$ws.Range("B5").Value = ''
$ws.Range("C5").Value = ''
$ws.Range("D5").Value = ''
$ws.Range("E5").Value = ''

$ws.Range("A6").Value = 'FC-DATE'
$ws.Range("B6").Value = ''
$ws.Range("C6").Value = ''
$ws.Range("D6").Value = ''
$ws.Range("E6").Value = ''

$ws.Range("A7").Value = 'FC-REV'
$ws.Range("B7").Value = ''
$ws.Range("C7").Value = ''
$ws.Range("D7").Value = ''
$ws.Range("E7").Value = ''

$ws.Range("A8").Value = 'FC-SC'
Set-TextValue $ws.Range("B8") '10'
$ws.Range("C8").Value = ''
$ws.Range("D8").Value = ''
$ws.Range("E8").Value = ''

$ws.Range("A9").Value = 'FC-SH'
Set-TextValue $ws.Range("B9") '1'
$ws.Range("C9").Value = 'X'
$ws.Range("D9").Value = ''
$ws.Range("E9").Value = ''

$ws.Range("A10").Value = 'FC-SI'
$ws.Range("B10").Value = 'A4'
$ws.Range("C10").Value = ''
$ws.Range("D10").Value = ''
$ws.Range("E10").Value = ''

$ws.Range("A11").Value = 'FreeCAD_DRAWING'
$ws.Range("B11").Value = 'FreeCAD DRAWING'
$ws.Range("C11").Value = ''
$ws.Range("D11").Value = ''
$ws.Range("E11").Value = ''

$ws.Range("A12").Value = 'Length_Units'
$ws.Range("B12").Value = 'mm'
$ws.Range("C12").Value = ''
$ws.Range("D12").Value = ''
$ws.Range("E12").Value = ''

$ws.Range("A13").Value = 'Mass_Units'
$ws.Range("B13").Value = 'kg'
$ws.Range("C13").Value = ''
$ws.Range("D13").Value = ''
$ws.Range("E13").Value = ''

$ws.Range("A14").Value = 'Number of sheets'
Set-TextValue $ws.Range("B14") '4'
$ws.Range("C14").Value = ''
$ws.Range("D14").Value = ''
$ws.Range("E14").Value = ''

$ws.Range("A15").Value = 'PN'
$ws.Range("B15").Value = 'PN'
$ws.Range("C15").Value = ''
$ws.Range("D15").Value = ''
$ws.Range("E15").Value = ''

$ws.Range("A16").Value = 'SI-1'
$ws.Range("B16").Value = ''
$ws.Range("C16").Value = ''
$ws.Range("D16").Value = ''
$ws.Range("E16").Value = ''

$ws.Range("A17").Value = 'SI-3'
$ws.Range("B17").Value = ''
$ws.Range("C17").Value = ''
$ws.Range("D17").Value = ''
$ws.Range("E17").Value = ''

$ws.Range("A18").Value = ''
$ws.Range("B18").Value = ''
$ws.Range("C18").Value = ''
$ws.Range("D18").Value = ''
$ws.Range("E18").Value = ''

# Grow the table (and its autofilter) so it covers the four extra rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E17"))

# Widen column A. The engine's ColumnWidth property bakes in the usual
# ~0.8333 "max digit width" padding offset, so asking for an on-disk
# width of 21 means setting ColumnWidth to 21 - 5/6.
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668
